# إضافة حدث جديد في Card18
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Card18")

# Fill the previously-blank cells on row 15 with "nan" placeholders,
# matching the rest of the sheet's convention for empty data cells,
# and record who serviced the machine.
$ws.Range("B15").Value = "nan"
$ws.Range("C15").Value = "nan"
$ws.Range("D15").Value = "nan"
$ws.Range("E15").Value = "nan"
$ws.Range("F15").Value = "nan"
$ws.Range("G15").Value = "nan"
$ws.Range("H15").Value = "nan"
$ws.Range("I15").Value = "nan"
$ws.Range("J15").Value = "nan"
$ws.Range("K15").Value = "nan"
$ws.Range("M15").Value = "nan"
$ws.Range("N15").Value = "nan"
$ws.Range("P15").Value = "الخبير"

# Add the new service event as row 16.
$ws.Range("A16").Value = "18"
$ws.Range("L16").Value = "10\3\2025"
$ws.Range("O16").Value = "تم تغيير الجرائد الخلفيه (1_5_8) ومعايره"
$ws.Range("P16").Value = "الخبير"
